# Update "diagramme de classes": add a new "Messagerie" class box
# (title + attributes) mirrored from the existing "Serveur" box, move
# the attributes that used to live on the Serveur box's second
# rectangle into the new Messagerie box, enlarge/empty the old
# rectangle and shift the methods rectangle down to make room.
#
# NB: Shape.Left/Top/Width/Height round-trip through a single-precision
# (float32) point value before being converted back to EMU for storage,
# truncating towards zero. The literals below are nudged by a few
# ULPs so that after that round-trip they land exactly on the target
# EMU values from the authoritative OOXML (still "the same" position,
# just compensating for the float32 storage the COM object model uses).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the existing shapes we need by name.
$attrsShape  = $null   # "Rectangle 6" (id 7)  - Serveur attributes box
$methodsShape = $null  # "Rectangle 8" (id 9)  - Serveur methods box
$titleShape  = $null   # "Rectangle 3" (id 4)  - Serveur title box

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Rectangle 6") { $attrsShape = $sh }
    if ($sh.Name -eq "Rectangle 8") { $methodsShape = $sh }
    if ($sh.Name -eq "Rectangle 3") { $titleShape = $sh }
}

# ---------------------------------------------------------------
# 1. New "Messagerie" title box, cloned from the "Serveur" title.
# ---------------------------------------------------------------
$newTitle = $titleShape.Duplicate().Item(1)
$newTitle.Name = "Rectangle 53"
$newTitle.Left = 377.2671661543307
$newTitle.Top = 23.117795275590552
$newTitle.Width = 158.2884292968504
$newTitle.Height = 25.17897707795276
$newTitle.TextFrame.TextRange.Text = "Messagerie"

# ---------------------------------------------------------------
# 2. New Messagerie attributes box, cloned from the old Serveur
#    attributes box (before it gets cleared out below) so the style
#    and run formatting match exactly.
# ---------------------------------------------------------------
$newAttrs = $attrsShape.Duplicate().Item(1)
$newAttrs.Name = "Rectangle 54"
$newAttrs.Left = 377.26708661417325
$newAttrs.Top = 48.29677165354331
$newAttrs.Width = 158.2885056370079
$newAttrs.Height = 45.592126884251975
$newAttrs.TextFrame.TextRange.Text = "# listeUsers : liste<String>`r# tableauId : HashMap<Integer, String>`r# tableauMessages : HashMap<Integer, String, int>"

# ---------------------------------------------------------------
# 3. New (empty, for now) Messagerie methods box, cloned again from
#    the Serveur attributes box, then text cleared.
# ---------------------------------------------------------------
$newMethods = $attrsShape.Duplicate().Item(1)
$newMethods.Name = "Rectangle 55"
$newMethods.Left = 377.26708661417325
$newMethods.Top = 93.88889763779528
$newMethods.Width = 158.2885056370079
$newMethods.Height = 47.1944103488189
$newMethods.TextFrame.TextRange.Text = ""

# ---------------------------------------------------------------
# 4. Resize/clear the original Serveur attributes box - its content
#    moved into the new Messagerie attributes box above.
# ---------------------------------------------------------------
$attrsShape.Top = 47.96007874015748
$attrsShape.Height = 50.9288196976378
$attrsShape.TextFrame.TextRange.Text = ""

# ---------------------------------------------------------------
# 5. Shift the Serveur methods box down to make room for the taller
#    attributes box above it.
# ---------------------------------------------------------------
$methodsShape.Top = 98.88889763779528
